# CCC19 Derived Variables Spreadsheet - add new ordinal outcome variable (O22a)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1. Insert a new blank row right after O22 (row 130), shifting O23 and subsequent rows down
$ws.Rows.Item(131).Insert()

# 2. Populate the new row with the new variable's data
$ws.Range("A131").Value = "O22a"
$ws.Range("B131").Value = "ordinal_v1a"
$ws.Range("C131").Value = "Outcome"
$ws.Range("D131").Value = "Custom ordinal outcome with death at any time"

# 3. Update the existing O22 row's description to reflect the "within 30 days" variant
$ws.Range("D130").Value = "Custom ordinal outcome with death within 30 days"

# 4. Expand the table (Table1) to include the newly inserted row
$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1 + 1
$tbl.Resize($ws.Range("A1:E" + $lastRow))

# 5. Leave the selection on the new row's Description cell, matching where the edit was made
$null = $ws.Range("D131").Select()
